$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.401.12'
$ws.Range("E2").Value = '  -5.72%  '
$ws.Range("D3").Value = '1.633.88'
$ws.Range("E3").Value = '  -7.22%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '304.89'
$ws.Range("E6").Value = '  -3.83%  '
$ws.Range("D7").Value = '0.3610'
$ws.Range("E7").Value = '  -5.67%  '
$ws.Range("D8").Value = '46.88'
$ws.Range("E8").Value = '  -6.96%  '
$ws.Range("D9").Value = '0.3217'
$ws.Range("E9").Value = '  -10.75%  '
$ws.Range("D10").Value = '1.100'
$ws.Range("E10").Value = '  -10.59%  '
$ws.Range("D11").Value = '0.06879'
$ws.Range("E11").Value = '  -10.32%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '5.888'
$ws.Range("E13").Value = '  -9.08%  '
$ws.Range("D14").Value = '19.08'
$ws.Range("E14").Value = '  -12.13%  '
$ws.Range("D15").Value = '1.634.59'
$ws.Range("E15").Value = '  -7.34%  '
$ws.Range("D16").Value = '6.498'
$ws.Range("E16").Value = '  -8.22%  '
$ws.Range("E17").Value = '  -10.03%  '
$ws.Range("D18").Value = '0.06506'
$ws.Range("E18").Value = '  -4.13%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '76.36'
$ws.Range("E20").Value = '  -12.49%  '
$ws.Range("D21").Value = '15.69'
$ws.Range("E21").Value = '  -11.26%  '
$ws.Range("D22").Value = '5.850'
$ws.Range("E22").Value = '  -10.03%  '
$ws.Range("D23").Value = '11.88'
$ws.Range("E23").Value = '  -7.95%  '
$ws.Range("D24").Value = '24.371.24'
$ws.Range("E24").Value = '  -5.54%  '
$ws.Range("D25").Value = '2.410'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").Value = '2.369'
$ws.Range("E26").Value = '  -18.62%  '
$ws.Range("D27").Value = '143.64'
$ws.Range("E27").Value = '  -7.87%  '
$ws.Range("D28").Value = '18.60'
$ws.Range("E28").Value = '  -10.11%  '
$ws.Range("D29").Value = '1.820.65'
$ws.Range("E29").Value = '  -7.17%  '
$ws.Range("E30").Value = '  -7.14%  '
$ws.Range("D31").Value = '1.089'
$ws.Range("E31").Value = '  -10.27%  '
$ws.Range("D32").Value = '4.063'
$ws.Range("E32").Value = '  -3.64%  '
$ws.Range("D33").Value = '5.619'
$ws.Range("E33").Value = '  -21.27%  '
$ws.Range("D34").Value = '0.08386'
$ws.Range("E34").Value = '  -4.20%  '
$ws.Range("E35").Value = '  -7.63%  '
$ws.Range("D36").Value = '12.27'
$ws.Range("E36").Value = '  -13.71%  '
$ws.Range("D37").Value = '5.096'
$ws.Range("E37").Value = '  -10.66%  '
$ws.Range("D38").Value = '0.05986'
$ws.Range("E38").Value = '  -11.34%  '
$ws.Range("D39").Value = '0.02202'
$ws.Range("E39").Value = '  -11.64%  '
$ws.Range("D40").Value = '1.196'
$ws.Range("E40").Value = '  -7.49%  '
$ws.Range("D41").Value = '0.2029'
$ws.Range("E41").Value = '  -10.09%  '
$ws.Range("D42").Value = '8.104'
$ws.Range("E42").Value = '  -13.51%  '
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '0.5820'
$ws.Range("E44").Value = '  -11.52%  '
$ws.Range("D45").Value = '3.715'
$ws.Range("E45").Value = '  -4.80%  '
$ws.Range("D46").Value = '12.47'
$ws.Range("E46").Value = '  -13.31%  '
$ws.Range("D47").Value = '0.5508'
$ws.Range("E47").Value = '  -13.20%  '
$ws.Range("D48").Value = '121.11'
$ws.Range("E48").Value = '  -8.44%  '
$ws.Range("D49").Value = '1.906'
$ws.Range("E49").Value = '  -12.18%  '
$ws.Range("D50").Value = '0.06898'
$ws.Range("E50").Value = '  -8.14%  '
$ws.Range("D51").Value = '73.16'
$ws.Range("E51").Value = '  -9.60%  '
